$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.954.55'
$ws.Range('E2').Value = '  +6.18%  '

$ws.Range('D3').Value = '3.545.58'
$ws.Range('E3').Value = '  +9.63%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '566.15'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +7.03%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '189.41'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +10.73%  '

$ws.Range('D7').Value = '3.538.51'
$ws.Range('E7').Value = '  +9.28%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.618'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +3.98%  '

$ws.Range('E9').Value = '  +0.03%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.634'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.82%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.151'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +13.73%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '54.80'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +3.32%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000270'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +6.27%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '9.42'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.40%  '

$ws.Range('D15').Value = '4.101.39'
$ws.Range('E15').Value = '  +9.51%  '

$ws.Range('D16').Value = '3.543.56'
$ws.Range('E16').Value = '  +9.58%  '

$ws.Range('E17').Value = '  +4.25%  '

$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '66.999.89'
$ws.Range('E18').Value = '  +6.49%  '

$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '18.29'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +6.40%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.04'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +8.84%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.999'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.64%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '434.02'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +18.48%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.16'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +11.17%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '85.23'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +5.48%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '4.13'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.43%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.12'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.06%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.89'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +9.52%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '12.24'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +9.08%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.14'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +11.49%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '30.51'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +7.34%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '641.85'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.20%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.60'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.87%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '11.76'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +4.99%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.112'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +5.83%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '59.88'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +5.33%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '38.52'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +5.63%  '

$ws.Range('D37').Value = '0.0₃0813'
$ws.Range('E37').Value = '  +12.34%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.147'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +18.31%  '

$ws.Range('E39').Value = '  -0.07%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.391'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +4.15%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.37'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +15.02%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').ClearFormats()

$ws.Range('D43').Value = '3.036.48'
$ws.Range('E43').Value = '  +5.82%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.66'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +4.56%  '

$ws.Range('E45').Value = '  +11.79%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.39'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +9.42%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0419'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +6.57%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.77'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +3.06%  '

$ws.Range('E49').Value = '  +5.81%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '143.98'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +7.81%  '

$ws.Range('E51').Value = '  +11.28%  '
